function Set-TextValue {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row reorder block (rows 46-48): Aptos/TheSandbox/BabyDogeCoin shuffle ---
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws "D46" "0.00000000118"
$ws.Range("E46").Value = "  +1.02%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws "D47" "7.045"
$ws.Range("E47").Value = "  -3.08%  "

$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws "D48" "0.4022"
$ws.Range("E48").Value = "  -0.15%  "

# --- Remaining price/volume updates ---
$ws.Range("D2").Value = "29.393.39"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.847.69"
$ws.Range("E3").Value = "  -0.11%  "
Set-TextValue $ws "D4" "0.9987"
$ws.Range("E4").Value = "  -0.03%  "
Set-TextValue $ws "D5" "240.41"
$ws.Range("E5").Value = "  -0.58%  "
Set-TextValue $ws "D6" "0.6300"
$ws.Range("E6").Value = "  +0.60%  "
Set-TextValue $ws "D7" "1.000"
$ws.Range("E7").Value = "  +0.01%  "
Set-TextValue $ws "D8" "0.07537"
$ws.Range("E8").Value = "  -0.05%  "
Set-TextValue $ws "D9" "0.2956"
$ws.Range("E9").Value = "  -0.66%  "
Set-TextValue $ws "D10" "24.46"
$ws.Range("E10").Value = "  +0.72%  "
Set-TextValue $ws "D11" "0.07719"
$ws.Range("D12").Value = "1.853.25"
$ws.Range("E12").Value = "  -2.27%  "
Set-TextValue $ws "D13" "4.992"
$ws.Range("E13").Value = "  -0.21%  "
Set-TextValue $ws "D14" "0.6843"
$ws.Range("E14").Value = "  -0.27%  "
Set-TextValue $ws "D15" "0.00001000"
$ws.Range("E15").Value = "  +2.20%  "
Set-TextValue $ws "D16" "82.94"
$ws.Range("E16").Value = "  -1.10%  "
Set-TextValue $ws "D17" "6.144"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "29.444.42"
$ws.Range("E18").Value = "  -0.37%  "
Set-TextValue $ws "D19" "228.49"
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("E20").Value = "  -0.34%  "
Set-TextValue $ws "D21" "0.9998"
Set-TextValue $ws "D22" "7.540"
$ws.Range("E22").Value = "  -0.75%  "
Set-TextValue $ws "D23" "1.0000"
$ws.Range("E23").Value = "  -0.01%  "
Set-TextValue $ws "D24" "156.97"
$ws.Range("E24").Value = "  +0.75%  "
Set-TextValue $ws "D25" "0.1397"
$ws.Range("E25").Value = "  +0.32%  "
Set-TextValue $ws "D26" "8.375"
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("E27").Value = "  -0.30%  "
Set-TextValue $ws "D28" "1.468"
$ws.Range("E28").Value = "  -0.86%  "
Set-TextValue $ws "D29" "0.05699"
$ws.Range("E29").Value = "  -2.43%  "
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("E31").Value = "  +0.64%  "
Set-TextValue $ws "D32" "4.017"
$ws.Range("E32").Value = "  -0.02%  "
Set-TextValue $ws "D33" "1.844"
Set-TextValue $ws "D34" "1.155"
$ws.Range("E34").Value = "  -1.28%  "
Set-TextValue $ws "D35" "0.7146"
$ws.Range("E35").Value = "  -0.51%  "
Set-TextValue $ws "D36" "2.587"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "1.254.19"
$ws.Range("E37").Value = "  +1.40%  "
Set-TextValue $ws "D38" "0.01813"
$ws.Range("E38").Value = "  +1.95%  "
Set-TextValue $ws "D39" "2.786"
$ws.Range("E39").Value = "  -0.41%  "
Set-TextValue $ws "D40" "0.9132"
$ws.Range("E40").Value = "  +0.15%  "
Set-TextValue $ws "D41" "6.190"
$ws.Range("E41").Value = "  +1.17%  "
Set-TextValue $ws "D42" "1.000"
$ws.Range("D43").Value = "2.024.06"
$ws.Range("E43").Value = "  -2.03%  "
Set-TextValue $ws "D44" "101.06"
$ws.Range("E44").Value = "  -0.56%  "
Set-TextValue $ws "D45" "66.15"
$ws.Range("E45").Value = "  -1.61%  "
Set-TextValue $ws "D49" "9.094"
$ws.Range("E49").Value = "  -0.56%  "
Set-TextValue $ws "D50" "1.690"
$ws.Range("E50").Value = "  -0.85%  "
Set-TextValue $ws "D51" "0.1128"
$ws.Range("E51").Value = "  +1.04%  "
